# 自动更新价格数据：在第2行插入最新一天的数据，原有数据整体下移一行。
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 插入新行（第2行），下方所有行自动下移一行，维度从 D59 扩展为 D60。
$ws.Rows.Item(2).Insert()

# 新插入的行会继承表头（第1行）的加粗/边框样式，清除掉让它和其它数据行一样没有样式。
$ws.Range("B2:D2").ClearFormats()

# A列需要保持纯文本日期（而不是被 Excel 自动识别成日期序列号），
# 所以先把该单元格数字格式设为文本，再写入日期字符串，写入后清除格式痕迹。
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-01-18"
$ws.Range("A2").ClearFormats()

# 其余三列（铁矿石 / 焦煤 / H型钢材）沿用最近一天同样的数值。
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
